# Auto-generated: updates cryptos list values to match the target commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.349.29"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "'3.549.30"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'618.30"
$ws.Range("E5").Value = "  +5.67%  "
$ws.Range("D6").Value = "'187.74"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'0.639"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").Value = "'53.83"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "'0.0000309"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "'9.71"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "'4.118.84"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "'616.60"
$ws.Range("E15").Value = "  +7.46%  "
$ws.Range("D16").Value = "'70.410.91"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'12.83"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("D18").Value = "'19.14"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'3.563.00"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "'17.70"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").Value = "'105.13"
$ws.Range("E23").Value = "  +10.84%  "
$ws.Range("D24").Value = "'4.73"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'5.09"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'3.03"
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("D27").Value = "'10.99"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("E28").Value = "  +9.08%  "
$ws.Range("D29").Value = "'34.43"
$ws.Range("E29").Value = "  +6.52%  "
$ws.Range("D30").Value = "'7.09"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").Value = "'12.54"
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").Value = "'64.27"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "'3.72"
$ws.Range("E34").Value = "  +16.17%  "
$ws.Range("D35").Value = "'3.17"
$ws.Range("E35").Value = "  -4.94%  "
$ws.Range("D36").Value = "'537.13"
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.402"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").Value = "'37.30"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'0.0₃0785"
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'3.547.81"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.57"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("D43").Value = "'0.139"
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("D45").Value = "'2.96"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +5.37%  "
$ws.Range("D47").Value = "'3.39"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").Value = "'9.06"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").Value = "'134.03"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").Value = "  -4.57%  "
